$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.622.14"
$ws.Range("E2").Value = "  -6.09%  "
$ws.Range("D3").Value = "2.218.15"
$ws.Range("E3").Value = "  -6.28%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'314.99"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").Value = "'98.48"
$ws.Range("E6").Value = "  -10.94%  "
$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = "  -8.92%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -9.25%  "
$ws.Range("D10").Value = "'36.37"
$ws.Range("E10").Value = "  -12.35%  "
$ws.Range("D11").Value = "'53.74"
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("D12").Value = "'0.0821"
$ws.Range("E12").Value = "  -10.80%  "
$ws.Range("D13").Value = "'7.66"
$ws.Range("E13").Value = "  -10.12%  "
$ws.Range("E14").Value = "  -3.88%  "
$ws.Range("D15").Value = "2.555.75"
$ws.Range("E15").Value = "  -6.32%  "
$ws.Range("E16").Value = "  -13.59%  "
$ws.Range("D17").Value = "'13.99"
$ws.Range("E17").Value = "  -9.45%  "
$ws.Range("D18").Value = "2.221.62"
$ws.Range("E18").Value = "  -6.13%  "
$ws.Range("D19").Value = "42.506.76"
$ws.Range("E19").Value = "  -6.35%  "
$ws.Range("D20").Value = "'13.73"
$ws.Range("E20").Value = "  +4.72%  "
$ws.Range("E21").Value = "  -10.13%  "
$ws.Range("D22").Value = "'6.50"
$ws.Range("E22").Value = "  -11.22%  "
$ws.Range("D23").Value = "'65.33"
$ws.Range("E23").Value = "  -11.20%  "
$ws.Range("D24").Value = "'3.17"
$ws.Range("E24").Value = "  -7.87%  "
$ws.Range("D25").Value = "'235.58"
$ws.Range("E25").Value = "  -9.57%  "
$ws.Range("D26").Value = "'2.13"
$ws.Range("E26").Value = "  -6.68%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  -8.95%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -7.03%  "
$ws.Range("D30").Value = "'6.43"
$ws.Range("E30").Value = "  -13.03%  "
$ws.Range("D31").Value = "'20.48"
$ws.Range("E31").Value = "  -8.81%  "
$ws.Range("D32").Value = "'0.0881"
$ws.Range("E32").Value = "  -9.07%  "
$ws.Range("D33").Value = "'157.86"
$ws.Range("E33").Value = "  -7.20%  "
$ws.Range("D34").Value = "'33.82"
$ws.Range("E34").Value = "  -11.03%  "
$ws.Range("D35").Value = "'2.74"
$ws.Range("E35").Value = "  -6.79%  "
$ws.Range("E36").Value = "  +7.53%  "
$ws.Range("E37").Value = "  -7.08%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.44"
$ws.Range("E38").Value = "  -8.11%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.85"
$ws.Range("E39").Value = "  +7.23%  "
$ws.Range("E40").Value = "  -11.72%  "
$ws.Range("D41").Value = "'3.54"
$ws.Range("E41").Value = "  -9.36%  "
$ws.Range("D42").Value = "'0.0320"
$ws.Range("E42").Value = "  -10.48%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "1.786.55"
$ws.Range("E44").Value = "  +9.78%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'11.97"
$ws.Range("E45").Value = "  -8.79%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "'87.77"
$ws.Range("E46").Value = "  -13.05%  "
$ws.Range("E47").Value = "  -11.69%  "
$ws.Range("D48").Value = "'77.79"
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("D49").Value = "'5.34"
$ws.Range("D50").Value = "'60.29"
$ws.Range("E50").Value = "  -14.05%  "
$ws.Range("D51").Value = "'15.85"
$ws.Range("E51").Value = "  +58.38%  "
